# Update faturamento anual data for the 2025 row (row 9) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 3039907.85
$ws.Range("C9").Value = 482006.99
$ws.Range("D9").Value = 3521914.84
$ws.Range("E9").Value = 13.68593540438928
$ws.Range("F9").Value = 86.31406459561073
$ws.Range("G9").Value = -53.41642781415277
$ws.Range("H9").Value = -45.1034661100548
$ws.Range("I9").Value = 30464
$ws.Range("J9").Value = 1297
$ws.Range("K9").Value = 31761
$ws.Range("L9").Value = 21892
$ws.Range("M9").Value = 160.8767970034716
$ws.Range("N9").Value = 9.833866009339731
